$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '264.24'
Set-TextValue $ws.Range('G2') '2'
Set-TextValue $ws.Range('D3') '22.69'
Set-TextValue $ws.Range('G3') '2'
Set-TextValue $ws.Range('D4') '6.269'
Set-TextValue $ws.Range('G4') '2'
Set-TextValue $ws.Range('D5') '0.06135'
Set-TextValue $ws.Range('G5') '2'
Set-TextValue $ws.Range('D6') '3.588'
Set-TextValue $ws.Range('G6') '2'
Set-TextValue $ws.Range('D7') '6.690'
Set-TextValue $ws.Range('G7') '2'
Set-TextValue $ws.Range('D8') '1.343'
Set-TextValue $ws.Range('G8') '2'
Set-TextValue $ws.Range('D9') '0.8311'
Set-TextValue $ws.Range('G9') '2'
Set-TextValue $ws.Range('D10') '0.01349'
Set-TextValue $ws.Range('G10') '2'
Set-TextValue $ws.Range('D11') '0.1582'
Set-TextValue $ws.Range('G11') '2'
Set-TextValue $ws.Range('D12') '0.08123'
Set-TextValue $ws.Range('G12') '2'
Set-TextValue $ws.Range('D13') '0.03365'
Set-TextValue $ws.Range('G13') '2'
Set-TextValue $ws.Range('D14') '0.03178'
Set-TextValue $ws.Range('G14') '2'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D15') '3.952'
$ws.Range('E15').Value = '14MCDexMCB'
Set-TextValue $ws.Range('G15') '2'
$ws.Range('B16').Value = 'BitMartToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D16') '0.09245'
$ws.Range('E16').Value = '15BitMartTokenBMX'
Set-TextValue $ws.Range('G16') '2'
Set-TextValue $ws.Range('D17') '0.001702'
Set-TextValue $ws.Range('G17') '2'
Set-TextValue $ws.Range('D18') '0.04870'
Set-TextValue $ws.Range('G18') '2'
Set-TextValue $ws.Range('D19') '0.006263'
Set-TextValue $ws.Range('G19') '2'
Set-TextValue $ws.Range('D20') '0.005271'
$ws.Range('E20').Value = '19HotbitTokenHTB'
Set-TextValue $ws.Range('G20') '2'
Set-TextValue $ws.Range('G21') '2'
Set-TextValue $ws.Range('D22') '0.0001504'
Set-TextValue $ws.Range('G22') '2'
Set-TextValue $ws.Range('D23') '3.778'
Set-TextValue $ws.Range('G23') '2'
Set-TextValue $ws.Range('D24') '2.320'
Set-TextValue $ws.Range('G24') '2'
Set-TextValue $ws.Range('D25') '0.3343'
Set-TextValue $ws.Range('G25') '2'
Set-TextValue $ws.Range('D26') '0.1243'
Set-TextValue $ws.Range('G26') '2'
Set-TextValue $ws.Range('D27') '0.0002688'
Set-TextValue $ws.Range('G27') '2'
Set-TextValue $ws.Range('G28') '2'
Set-TextValue $ws.Range('G29') '2'
Set-TextValue $ws.Range('G30') '2'
Set-TextValue $ws.Range('G31') '2'
Set-TextValue $ws.Range('G32') '2'
Set-TextValue $ws.Range('G33') '2'
Set-TextValue $ws.Range('G34') '2'
Set-TextValue $ws.Range('G35') '2'
Set-TextValue $ws.Range('G36') '2'
Set-TextValue $ws.Range('G37') '2'
Set-TextValue $ws.Range('G38') '2'
Set-TextValue $ws.Range('G39') '2'
Set-TextValue $ws.Range('D40') '0.04617'
Set-TextValue $ws.Range('G40') '2'
Set-TextValue $ws.Range('D41') '0.006979'
Set-TextValue $ws.Range('G41') '2'
Set-TextValue $ws.Range('D42') '0.1134'
Set-TextValue $ws.Range('G42') '2'
Set-TextValue $ws.Range('D43') '0.003379'
Set-TextValue $ws.Range('G43') '2'
Set-TextValue $ws.Range('D44') '0.01215'
Set-TextValue $ws.Range('G44') '2'
Set-TextValue $ws.Range('D45') '0.00006152'
Set-TextValue $ws.Range('G45') '2'
Set-TextValue $ws.Range('D46') '0.00000000752'
Set-TextValue $ws.Range('G46') '2'
Set-TextValue $ws.Range('D47') '0.7914'
Set-TextValue $ws.Range('G47') '2'
Set-TextValue $ws.Range('D48') '0.1912'
Set-TextValue $ws.Range('G48') '2'
Set-TextValue $ws.Range('D49') '0.00001404'
$ws.Range('E49').Value = '48CryptobidCoinCBCWorstin24h'
Set-TextValue $ws.Range('G49') '2'
Set-TextValue $ws.Range('D50') '0.01244'
Set-TextValue $ws.Range('G50') '2'
Set-TextValue $ws.Range('G51') '2'
